$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = 0.5050505050505051
$ws.Cells.Item(2, 6).Value = 0.5056116722783389
$ws.Cells.Item(2, 7).Value = 0.5050505050505051
$ws.Cells.Item(2, 8).Value = 0.5040388555973635

$ws.Cells.Item(3, 5).Value = 0.5454545454545454
$ws.Cells.Item(3, 6).Value = 0.5532598714416896
$ws.Cells.Item(3, 7).Value = 0.5454545454545454
$ws.Cells.Item(3, 8).Value = 0.5325015476927473

$ws.Cells.Item(4, 5).Value = 0.5757575757575758
$ws.Cells.Item(4, 6).Value = 0.5760546642899584
$ws.Cells.Item(4, 7).Value = 0.5757575757575758
$ws.Cells.Item(4, 8).Value = 0.5756709956709957

$ws.Cells.Item(5, 5).Value = 0.6565656565656566
$ws.Cells.Item(5, 6).Value = 0.6569617746088333
$ws.Cells.Item(5, 7).Value = 0.6565656565656566
$ws.Cells.Item(5, 8).Value = 0.6564955679241393

$ws.Cells.Item(6, 5).Value = 0.5858585858585859
$ws.Cells.Item(6, 6).Value = 0.5882076579750998
$ws.Cells.Item(6, 7).Value = 0.5858585858585859
$ws.Cells.Item(6, 8).Value = 0.5840773324644293

$ws.Cells.Item(7, 5).Value = 0.6262626262626263
$ws.Cells.Item(7, 6).Value = 0.6228474368009251
$ws.Cells.Item(7, 7).Value = 0.6262626262626263
$ws.Cells.Item(7, 8).Value = 0.5658616089325262

$ws.Cells.Item(9, 5).Value = 0.5858585858585859
$ws.Cells.Item(9, 6).Value = 0.5644007644007645
$ws.Cells.Item(9, 7).Value = 0.5858585858585859
$ws.Cells.Item(9, 8).Value = 0.5614268772163509

$ws.Cells.Item(10, 5).Value = 0.6060606060606061
$ws.Cells.Item(10, 6).Value = 0.5886158886158885
$ws.Cells.Item(10, 7).Value = 0.6060606060606061
$ws.Cells.Item(10, 8).Value = 0.582820688083846

$ws.Cells.Item(11, 5).Value = 0.5858585858585859
$ws.Cells.Item(11, 6).Value = 0.5617906756350692
$ws.Cells.Item(11, 7).Value = 0.5858585858585859
$ws.Cells.Item(11, 8).Value = 0.5560579338357117

$ws.Cells.Item(12, 5).Value = 0.4747474747474748
$ws.Cells.Item(12, 6).Value = 0.4753494124922696
$ws.Cells.Item(12, 7).Value = 0.4747474747474748
$ws.Cells.Item(12, 8).Value = 0.4748546691403834

$ws.Cells.Item(13, 5).Value = 0.5757575757575758
$ws.Cells.Item(13, 6).Value = 0.5752066115702479
$ws.Cells.Item(13, 7).Value = 0.5757575757575758
$ws.Cells.Item(13, 8).Value = 0.5745394884033111

$ws.Cells.Item(14, 5).Value = 0.6363636363636364
$ws.Cells.Item(14, 6).Value = 0.6398933249666641
$ws.Cells.Item(14, 7).Value = 0.6363636363636364
$ws.Cells.Item(14, 8).Value = 0.6314557535487768

$ws.Cells.Item(15, 5).Value = 0.6363636363636364
$ws.Cells.Item(15, 6).Value = 0.6475524475524476
$ws.Cells.Item(15, 7).Value = 0.6363636363636364
$ws.Cells.Item(15, 8).Value = 0.6323232323232323

$ws.Cells.Item(16, 5).Value = 0.6262626262626263
$ws.Cells.Item(16, 6).Value = 0.6354453627180899
$ws.Cells.Item(16, 7).Value = 0.6262626262626263
$ws.Cells.Item(16, 8).Value = 0.6156148378370602

$ws.Cells.Item(17, 5).Value = 0.7676767676767676
$ws.Cells.Item(17, 6).Value = 0.7348484848484849
$ws.Cells.Item(17, 7).Value = 0.7676767676767676
$ws.Cells.Item(17, 8).Value = 0.7337954479058773

$ws.Cells.Item(18, 6).Value = 0.722048066875653
$ws.Cells.Item(18, 8).Value = 0.7261503928170594

$ws.Cells.Item(19, 5).Value = 0.6464646464646465
$ws.Cells.Item(19, 6).Value = 0.6063432024694578
$ws.Cells.Item(19, 7).Value = 0.6464646464646465
$ws.Cells.Item(19, 8).Value = 0.6241661605961747

$ws.Cells.Item(24, 5).Value = 0.6363636363636364
$ws.Cells.Item(24, 6).Value = 0.6394219741570457
$ws.Cells.Item(24, 7).Value = 0.6363636363636364
$ws.Cells.Item(24, 8).Value = 0.6336700336700336

$ws.Cells.Item(25, 5).Value = 0.5151515151515151
$ws.Cells.Item(25, 6).Value = 0.5408432147562582
$ws.Cells.Item(25, 7).Value = 0.5151515151515151
$ws.Cells.Item(25, 8).Value = 0.4050362782757149

$ws.Cells.Item(26, 5).Value = 0.5656565656565656
$ws.Cells.Item(26, 6).Value = 0.5813243073517046
$ws.Cells.Item(26, 7).Value = 0.5656565656565656
$ws.Cells.Item(26, 8).Value = 0.5396665845446332

$ws.Cells.Item(27, 5).Value = 0.5959595959595959
$ws.Cells.Item(27, 6).Value = 0.6085378673613967
$ws.Cells.Item(27, 7).Value = 0.5959595959595959
$ws.Cells.Item(27, 8).Value = 0.5858332700437965

$ws.Cells.Item(28, 5).Value = 0.6060606060606061
$ws.Cells.Item(28, 6).Value = 0.5845615408025738
$ws.Cells.Item(28, 7).Value = 0.6060606060606061
$ws.Cells.Item(28, 8).Value = 0.5587114142545274

$ws.Cells.Item(29, 5).Value = 0.6161616161616161
$ws.Cells.Item(29, 6).Value = 0.6173600410888547
$ws.Cells.Item(29, 7).Value = 0.6161616161616161
$ws.Cells.Item(29, 8).Value = 0.61267217630854

$ws.Cells.Item(30, 5).Value = 0.7474747474747475
$ws.Cells.Item(30, 6).Value = 0.6772404900064475
$ws.Cells.Item(30, 7).Value = 0.7474747474747475
$ws.Cells.Item(30, 8).Value = 0.678946164357305

$ws.Cells.Item(31, 6).Value = 0.6588991177598773
$ws.Cells.Item(31, 8).Value = 0.5675562500638335

